# Add a "Save" column (H) to the s_vals sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell, styled the same way as the other header cells (B1:G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Data values for the new "Save" column
$saveValues = @(0, 0, 0, 1, 0, 1, 0, 0, 0, 0, 0, 1)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
